$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from G1 into H1, then set its value
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Save"

# Fill the new "Save" column with zeros
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
